$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.063.42'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.30%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.645.69'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.68'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5215'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2610'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.62%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06359'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.81'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.69%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07661'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.645.06'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.43%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.421'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.868.57'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5541'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8252'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.95'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.082.20'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.24%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.731'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.38%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '188.63'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.95%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.19'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.231'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.15%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.16'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1219'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.420'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.84'

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05957'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.42%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.268'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.398'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.395'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.73%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.662'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.89%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9954'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.92%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.393'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.755'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5615'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -7.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01613'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.05%  '

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.841'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.15%  '

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8550'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.028.97'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -7.66%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.13'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.41%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.795.29'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.47%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈110'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.75'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.44%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9979'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.47%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.050'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05153'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.59%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4220'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.52%  '
